{"js": "// \"Miswrite the T number\": the DFA-node paragraph that was mistakenly\n// labelled \"T72:\" (right after the \"ELSE\" STMT rule / right before the\n// genuine \"T72:\" paragraph) must read \"T71:\" instead. Also the stray\n// \"_GoBack\" bookmark left over from the last edit position is removed.\n\n// 1) Remove the leftover \"_GoBack\" bookmark (now-empty paragraph stays,\n//    but its bookmarkStart/bookmarkEnd children go away).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Fix the mislabeled paragraph: the first \"T72:\" (the one following\n//    the \"if ... rbrace ELSE\" STMT paragraph) should actually be \"T71:\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === \"T72:\") {\n    const prev = para.getPreviousOrNullObject();\n    prev.load(\"text\");\n    await context.sync();\n    if (!prev.isNullObject) {\n      const prevPrev = prev.getPreviousOrNullObject();\n      prevPrev.load(\"text\");\n      await context.sync();\n      if (!prevPrev.isNullObject && prevPrev.text.indexOf(\"rbrace ELSE\") !== -1) {\n        target = para;\n        break;\n      }\n    }\n  }\n}\n\nif (target) {\n  const results = target.search(\"72:\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  const match = results.items[0];\n\n  // Replace the single \"72:\" run with three separate runs (\"7\", \"1\", \":\")\n  // exactly as the authored edit did, keeping the leading \"T\" run intact.\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships></pkg:xmlData></pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p><w:r><w:t>7</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p></w:body>' +\n    '</w:document></pkg:xmlData></pkg:part></pkg:package>';\n\n  match.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"Miswrite the T number\": the DFA-node paragraph mistakenly labelled\n# \"T72:\" (the one right after the \"if ... rbrace ELSE\" STMT rule, and\n# right before the genuine \"T72:\" paragraph) must actually read \"T71:\".\n# Also removes the stray \"_GoBack\" bookmark left over from the last\n# saved cursor position.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the leftover \"_GoBack\" bookmark (the paragraph that hosted\n#    it stays, now empty of any bookmark markup).\ntry {\n    $goBack = $d.Bookmarks.Item(\"_GoBack\")\n    $goBack.Delete()\n} catch {\n    # Bookmark already absent - nothing to do.\n}\n\n# 2) Locate the mislabeled \"T72:\" paragraph - the first one, which\n#    immediately follows the \"... rbrace ELSE\" STMT paragraph - and\n#    fix it to read \"T71:\".\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($ptext -ne \"T72:\") {\n        continue\n    }\n\n    $prev = $p.Previous()\n    $prevText = \"\"\n    if ($prev -ne $null) {\n        $prevText = $prev.Range.Text.TrimEnd([char]13, [char]7)\n    }\n    $prevPrevText = \"\"\n    if ($prev -ne $null -and $prevText -eq \"\") {\n        $prevPrev = $prev.Previous()\n        if ($prevPrev -ne $null) {\n            $prevPrevText = $prevPrev.Range.Text.TrimEnd([char]13, [char]7)\n        }\n    }\n\n    if ($prevPrevText -notlike \"*rbrace ELSE*\") {\n        continue\n    }\n\n    # Narrow a Range to just the \"72:\" run text within this paragraph.\n    $searchRange = $p.Range\n    $found = $searchRange.Find.Execute(\"72:\")\n    if (-not $found) {\n        continue\n    }\n\n    # Re-anchor a fresh Range over the same offsets (InsertXML replaces\n    # the exact contents of the Range it is called on).\n    $target = $d.Range($searchRange.Start, $searchRange.End)\n\n    $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>7</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n    $target.InsertXML($xml)\n    break\n}\n"}
